# Insert a new row for "HSD010 - General health condition" above the
# existing "Self-perception summary stat" row (old row 70), shifting the
# remaining rows (old 70-74) down by one (new 71-75).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 70; Excel copies formatting from the row
# above it (row 69), which already uses style index 2 (yellow fill) on
# columns B and C - matching the desired formatting for the new row.
$ws.Rows.Item(70).Insert() | Out-Null

# Populate the newly inserted row 70 with the new variable information.
$ws.Cells.Item(70, 1).Value = "Qns"
$ws.Cells.Item(70, 2).Value = "HSD010"
$ws.Cells.Item(70, 3).Value = "General health condition"

# Update the view to reflect where the edit was made.
$ws.Range("C71").Select() | Out-Null
